$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44757
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 25000
$ws.Cells.Item(2, 12).Value = 25000
$ws.Cells.Item(2, 13).Value = 25000
$ws.Cells.Item(2, 16).Value = 1667

# Row 3
$ws.Cells.Item(3, 4).Value = 44754
$ws.Cells.Item(3, 10).Value = 90

# Row 4
$ws.Cells.Item(4, 4).Value = 44792
$ws.Cells.Item(4, 10).Value = 120
$ws.Cells.Item(4, 11).Value = 24000
$ws.Cells.Item(4, 12).Value = 24000
$ws.Cells.Item(4, 13).Value = 24000
$ws.Cells.Item(4, 16).Value = 1600

# Row 5
$ws.Cells.Item(5, 4).Value = 44819
$ws.Cells.Item(5, 11).Value = 22000
$ws.Cells.Item(5, 12).Value = 22000
$ws.Cells.Item(5, 13).Value = 22000
$ws.Cells.Item(5, 16).Value = 1467

# Row 6
$ws.Cells.Item(6, 4).Value = 45163
$ws.Cells.Item(6, 10).Value = 140
$ws.Cells.Item(6, 11).Value = 25000
$ws.Cells.Item(6, 13).Value = 25000
$ws.Cells.Item(6, 16).Value = 1667

# Row 7
$ws.Cells.Item(7, 4).Value = 44803
$ws.Cells.Item(7, 10).Value = 90
$ws.Cells.Item(7, 11).Value = 24000
$ws.Cells.Item(7, 12).Value = 24000
$ws.Cells.Item(7, 13).Value = 24000
$ws.Cells.Item(7, 16).Value = 1600

# Row 8
$ws.Cells.Item(8, 4).Value = 44761
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 23000
$ws.Cells.Item(8, 12).Value = 25000

# Row 9
$ws.Cells.Item(9, 4).Value = 45167

# Row 10
$ws.Cells.Item(10, 4).Value = 45156
$ws.Cells.Item(10, 10).Value = 120
$ws.Cells.Item(10, 11).Value = 25000
$ws.Cells.Item(10, 12).Value = 25000
$ws.Cells.Item(10, 13).Value = 25000
$ws.Cells.Item(10, 16).Value = 1667

# Row 11
$ws.Cells.Item(11, 4).Value = 44782
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 24000
$ws.Cells.Item(11, 12).Value = 24000
$ws.Cells.Item(11, 13).Value = 24000
$ws.Cells.Item(11, 16).Value = 1600

# Row 12
$ws.Cells.Item(12, 4).Value = 44771
$ws.Cells.Item(12, 11).Value = 25000
$ws.Cells.Item(12, 12).Value = 25000
$ws.Cells.Item(12, 13).Value = 25000
$ws.Cells.Item(12, 16).Value = 1667

# Row 13
$ws.Cells.Item(13, 4).Value = 44827
$ws.Cells.Item(13, 10).Value = 90
$ws.Cells.Item(13, 11).Value = 22000
$ws.Cells.Item(13, 12).Value = 22000
$ws.Cells.Item(13, 13).Value = 22000
$ws.Cells.Item(13, 16).Value = 1467

# Row 14
$ws.Cells.Item(14, 4).Value = 44775
$ws.Cells.Item(14, 10).Value = 120
$ws.Cells.Item(14, 11).Value = 24000
$ws.Cells.Item(14, 12).Value = 24000
$ws.Cells.Item(14, 13).Value = 24000
$ws.Cells.Item(14, 16).Value = 1600

# Row 15
$ws.Cells.Item(15, 4).Value = 44831
$ws.Cells.Item(15, 10).Value = 90
$ws.Cells.Item(15, 11).Value = 25000
$ws.Cells.Item(15, 12).Value = 25000
$ws.Cells.Item(15, 13).Value = 25000
$ws.Cells.Item(15, 16).Value = 1667

# Row 16
$ws.Cells.Item(16, 4).Value = 44810
$ws.Cells.Item(16, 10).Value = 110
$ws.Cells.Item(16, 11).Value = 22000
$ws.Cells.Item(16, 12).Value = 22000
$ws.Cells.Item(16, 13).Value = 22000
$ws.Cells.Item(16, 16).Value = 1467

# Row 17
$ws.Cells.Item(17, 4).Value = 44838
$ws.Cells.Item(17, 10).Value = 80
$ws.Cells.Item(17, 11).Value = 22000
$ws.Cells.Item(17, 12).Value = 22000
$ws.Cells.Item(17, 13).Value = 22000
$ws.Cells.Item(17, 16).Value = 1467

# Row 18
$ws.Cells.Item(18, 4).Value = 44799
$ws.Cells.Item(18, 10).Value = 80
$ws.Cells.Item(18, 11).Value = 23000
$ws.Cells.Item(18, 12).Value = 23000
$ws.Cells.Item(18, 13).Value = 23000
$ws.Cells.Item(18, 16).Value = 1533

# Row 19
$ws.Cells.Item(19, 4).Value = 44407
$ws.Cells.Item(19, 10).Value = 90

# Row 20
$ws.Cells.Item(20, 4).Value = 44764
$ws.Cells.Item(20, 10).Value = 90
$ws.Cells.Item(20, 11).Value = 24000
$ws.Cells.Item(20, 12).Value = 24000
$ws.Cells.Item(20, 13).Value = 24000
$ws.Cells.Item(20, 16).Value = 1600

# Row 21
$ws.Cells.Item(21, 4).Value = 44817
$ws.Cells.Item(21, 10).Value = 90
$ws.Cells.Item(21, 11).Value = 23000
$ws.Cells.Item(21, 12).Value = 23000
$ws.Cells.Item(21, 13).Value = 23000
$ws.Cells.Item(21, 16).Value = 1533

# Row 22
$ws.Cells.Item(22, 4).Value = 44365
$ws.Cells.Item(22, 10).Value = 80
$ws.Cells.Item(22, 11).Value = 25000
$ws.Cells.Item(22, 12).Value = 25000
$ws.Cells.Item(22, 13).Value = 25000
$ws.Cells.Item(22, 16).Value = 1667

# Row 23
$ws.Cells.Item(23, 4).Value = 44750

# Row 24
$ws.Cells.Item(24, 4).Value = 45177
$ws.Cells.Item(24, 10).Value = 120
$ws.Cells.Item(24, 11).Value = 26000
$ws.Cells.Item(24, 12).Value = 26000
$ws.Cells.Item(24, 13).Value = 26000
$ws.Cells.Item(24, 16).Value = 1733

# Row 25
$ws.Cells.Item(25, 4).Value = 45146
$ws.Cells.Item(25, 10).Value = 140
$ws.Cells.Item(25, 11).Value = 26000
$ws.Cells.Item(25, 12).Value = 26000
$ws.Cells.Item(25, 13).Value = 26000
$ws.Cells.Item(25, 16).Value = 1733

# Row 26
$ws.Cells.Item(26, 4).Value = 44806
$ws.Cells.Item(26, 10).Value = 70
$ws.Cells.Item(26, 11).Value = 23000
$ws.Cells.Item(26, 12).Value = 23000
$ws.Cells.Item(26, 13).Value = 23000
$ws.Cells.Item(26, 16).Value = 1533

# Row 27
$ws.Cells.Item(27, 4).Value = 44418
$ws.Cells.Item(27, 10).Value = 90
$ws.Cells.Item(27, 11).Value = 25000
$ws.Cells.Item(27, 12).Value = 25000
$ws.Cells.Item(27, 13).Value = 25000
$ws.Cells.Item(27, 16).Value = 1667

# Row 28
$ws.Cells.Item(28, 4).Value = 44778
$ws.Cells.Item(28, 10).Value = 120
$ws.Cells.Item(28, 11).Value = 24000
$ws.Cells.Item(28, 12).Value = 24000
$ws.Cells.Item(28, 13).Value = 24000
$ws.Cells.Item(28, 16).Value = 1600

# Row 29
$ws.Cells.Item(29, 4).Value = 44789
$ws.Cells.Item(29, 10).Value = 90
$ws.Cells.Item(29, 11).Value = 24000
$ws.Cells.Item(29, 12).Value = 24000
$ws.Cells.Item(29, 13).Value = 24000
$ws.Cells.Item(29, 16).Value = 1600

# Row 30
$ws.Cells.Item(30, 4).Value = 44400
$ws.Cells.Item(30, 10).Value = 80

# Row 31
$ws.Cells.Item(31, 4).Value = 44740
$ws.Cells.Item(31, 10).Value = 90
$ws.Cells.Item(31, 11).Value = 25000
$ws.Cells.Item(31, 12).Value = 25000
$ws.Cells.Item(31, 13).Value = 25000
$ws.Cells.Item(31, 16).Value = 1667

# Row 32
$ws.Cells.Item(32, 4).Value = 44781
$ws.Cells.Item(32, 10).Value = 70
$ws.Cells.Item(32, 11).Value = 24000
$ws.Cells.Item(32, 12).Value = 24000
$ws.Cells.Item(32, 13).Value = 24000
$ws.Cells.Item(32, 16).Value = 1600

# Row 33
$ws.Cells.Item(33, 4).Value = 45149
$ws.Cells.Item(33, 10).Value = 120
